$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for F1:H1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style (bold, centered, bordered) from A1 onto the new header cells
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill the boolean "Outliers_MAD" columns for rows 2-12 (all False, except F12 = True)
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
$ws.Cells.Item(12, 6).Value = $true
